# Auto-generated data-driven edit list applying the diff to Moogle_Profits workbook.
# Each entry: worksheet name, A1 cell reference, and the new literal value
# (or $null to clear/remove the cell's content entirely, matching a cell
# deletion in the source diff).
$wb = $excel.ActiveWorkbook

$edits = @(
    @{ Sheet = "ALC"; Cell = "H8"; Value = 6.8333335 }
    @{ Sheet = "ALC"; Cell = "I8"; Value = 6.8333335 }
    @{ Sheet = "ALC"; Cell = "K8"; Value = 20.5000005 }
    @{ Sheet = "ALC"; Cell = "M8"; Value = 118.4999995 }
    @{ Sheet = "ALC"; Cell = "H11"; Value = 1.8 }
    @{ Sheet = "ALC"; Cell = "I11"; Value = 1.8 }
    @{ Sheet = "ALC"; Cell = "K11"; Value = 1.8 }
    @{ Sheet = "ALC"; Cell = "M11"; Value = 138.2 }
    @{ Sheet = "ALC"; Cell = "H92"; Value = 502.5 }
    @{ Sheet = "ALC"; Cell = "I92"; Value = 132.14285 }
    @{ Sheet = "ALC"; Cell = "J92"; Value = 1366.6666 }
    @{ Sheet = "ALC"; Cell = "K92"; Value = 132.14285 }
    @{ Sheet = "ALC"; Cell = "L92"; Value = 1366.6666 }
    @{ Sheet = "ALC"; Cell = "M92"; Value = 1115.85715 }
    @{ Sheet = "ALC"; Cell = "N92"; Value = -3862.6666 }
    @{ Sheet = "ALC"; Cell = "H107"; Value = 431.87097 }
    @{ Sheet = "ALC"; Cell = "I107"; Value = 406.44 }
    @{ Sheet = "ALC"; Cell = "J107"; Value = 537.8333 }
    @{ Sheet = "ALC"; Cell = "K107"; Value = 406.44 }
    @{ Sheet = "ALC"; Cell = "L107"; Value = 537.8333 }
    @{ Sheet = "ALC"; Cell = "M107"; Value = 1513.56 }
    @{ Sheet = "ALC"; Cell = "N107"; Value = -4377.8333 }
    @{ Sheet = "ALC"; Cell = "H129"; Value = 947.2308 }
    @{ Sheet = "ALC"; Cell = "I129"; Value = 776.1667 }
    @{ Sheet = "ALC"; Cell = "J129"; Value = 3000 }
    @{ Sheet = "ALC"; Cell = "K129"; Value = 2328.5001 }
    @{ Sheet = "ALC"; Cell = "L129"; Value = 9000 }
    @{ Sheet = "ALC"; Cell = "M129"; Value = 2671.4999 }
    @{ Sheet = "ALC"; Cell = "N129"; Value = -19000 }
    @{ Sheet = "ALC"; Cell = "H131"; Value = 1879.375 }
    @{ Sheet = "ALC"; Cell = "J131"; Value = 3925 }
    @{ Sheet = "ALC"; Cell = "L131"; Value = 11775 }
    @{ Sheet = "ALC"; Cell = "N131"; Value = -21855 }
    @{ Sheet = "ALC"; Cell = "H137"; Value = 1506.24 }
    @{ Sheet = "ALC"; Cell = "I137"; Value = 1403.2174 }
    @{ Sheet = "ALC"; Cell = "K137"; Value = 4209.6522 }
    @{ Sheet = "ALC"; Cell = "M137"; Value = -1659.6522 }
    @{ Sheet = "ARM"; Cell = "H32"; Value = 13080.873 }
    @{ Sheet = "ARM"; Cell = "I32"; Value = 9000.378000000001 }
    @{ Sheet = "ARM"; Cell = "K32"; Value = 9000.378000000001 }
    @{ Sheet = "ARM"; Cell = "M32"; Value = -8713.378000000001 }
    @{ Sheet = "ARM"; Cell = "H45"; Value = 3573839.2 }
    @{ Sheet = "ARM"; Cell = "I45"; Value = 4349700 }
    @{ Sheet = "ARM"; Cell = "K45"; Value = 4349700 }
    @{ Sheet = "ARM"; Cell = "M45"; Value = -4349323 }
    @{ Sheet = "ARM"; Cell = "H46"; Value = 7866 }
    @{ Sheet = "ARM"; Cell = "J46"; Value = 8849.200000000001 }
    @{ Sheet = "ARM"; Cell = "L46"; Value = 8849.200000000001 }
    @{ Sheet = "ARM"; Cell = "N46"; Value = -9487.200000000001 }
    @{ Sheet = "ARM"; Cell = "H61"; Value = 7933 }
    @{ Sheet = "ARM"; Cell = "I61"; Value = 7562.5713 }
    @{ Sheet = "ARM"; Cell = "K61"; Value = 7562.5713 }
    @{ Sheet = "ARM"; Cell = "M61"; Value = -7350.5713 }
    @{ Sheet = "ARM"; Cell = "H74"; Value = 12807.111 }
    @{ Sheet = "ARM"; Cell = "I74"; Value = 5233 }
    @{ Sheet = "ARM"; Cell = "J74"; Value = 16594.166 }
    @{ Sheet = "ARM"; Cell = "K74"; Value = 5233 }
    @{ Sheet = "ARM"; Cell = "L74"; Value = 16594.166 }
    @{ Sheet = "ARM"; Cell = "M74"; Value = -4359 }
    @{ Sheet = "ARM"; Cell = "N74"; Value = -18342.166 }
    @{ Sheet = "ARM"; Cell = "H77"; Value = 12807.111 }
    @{ Sheet = "ARM"; Cell = "I77"; Value = 5233 }
    @{ Sheet = "ARM"; Cell = "J77"; Value = 16594.166 }
    @{ Sheet = "ARM"; Cell = "K77"; Value = 26165 }
    @{ Sheet = "ARM"; Cell = "L77"; Value = 82970.83 }
    @{ Sheet = "ARM"; Cell = "M77"; Value = -21797 }
    @{ Sheet = "ARM"; Cell = "N77"; Value = -91706.83 }
    @{ Sheet = "ARM"; Cell = "H122"; Value = 2644.7144 }
    @{ Sheet = "ARM"; Cell = "I122"; Value = 2186.74 }
    @{ Sheet = "ARM"; Cell = "K122"; Value = 6560.219999999999 }
    @{ Sheet = "ARM"; Cell = "M122"; Value = -4110.219999999999 }
    @{ Sheet = "ARM"; Cell = "H132"; Value = 4471.2964 }
    @{ Sheet = "ARM"; Cell = "I132"; Value = 2486.9524 }
    @{ Sheet = "ARM"; Cell = "J132"; Value = 11416.5 }
    @{ Sheet = "ARM"; Cell = "K132"; Value = 7460.8572 }
    @{ Sheet = "ARM"; Cell = "L132"; Value = 34249.5 }
    @{ Sheet = "ARM"; Cell = "M132"; Value = -4930.8572 }
    @{ Sheet = "ARM"; Cell = "N132"; Value = -39309.5 }
    @{ Sheet = "ARM"; Cell = "H136"; Value = 7933 }
    @{ Sheet = "ARM"; Cell = "I136"; Value = 7562.5713 }
    @{ Sheet = "ARM"; Cell = "K136"; Value = 22687.7139 }
    @{ Sheet = "ARM"; Cell = "M136"; Value = -20137.7139 }
    @{ Sheet = "BSM"; Cell = "H86"; Value = 3705.52 }
    @{ Sheet = "BSM"; Cell = "I86"; Value = 2114.5881 }
    @{ Sheet = "BSM"; Cell = "K86"; Value = 2114.5881 }
    @{ Sheet = "BSM"; Cell = "M86"; Value = -991.5880999999999 }
    @{ Sheet = "BSM"; Cell = "H89"; Value = 3705.52 }
    @{ Sheet = "BSM"; Cell = "I89"; Value = 2114.5881 }
    @{ Sheet = "BSM"; Cell = "K89"; Value = 10572.9405 }
    @{ Sheet = "BSM"; Cell = "M89"; Value = -4956.940500000001 }
    @{ Sheet = "BSM"; Cell = "H99"; Value = 2514.4546 }
    @{ Sheet = "BSM"; Cell = "I99"; Value = 2606.4285 }
    @{ Sheet = "BSM"; Cell = "K99"; Value = 2606.4285 }
    @{ Sheet = "BSM"; Cell = "M99"; Value = -1108.4285 }
    @{ Sheet = "CRP"; Cell = "H31"; Value = 8467.102999999999 }
    @{ Sheet = "CRP"; Cell = "I31"; Value = 3149.5 }
    @{ Sheet = "CRP"; Cell = "K31"; Value = 3149.5 }
    @{ Sheet = "CRP"; Cell = "M31"; Value = -2854.5 }
    @{ Sheet = "CRP"; Cell = "H34"; Value = 8467.102999999999 }
    @{ Sheet = "CRP"; Cell = "I34"; Value = 3149.5 }
    @{ Sheet = "CRP"; Cell = "K34"; Value = 3149.5 }
    @{ Sheet = "CRP"; Cell = "M34"; Value = -2947.5 }
    @{ Sheet = "CRP"; Cell = "H86"; Value = 3082420 }
    @{ Sheet = "CRP"; Cell = "J86"; Value = 7059.4 }
    @{ Sheet = "CRP"; Cell = "L86"; Value = 7059.4 }
    @{ Sheet = "CRP"; Cell = "N86"; Value = -9305.4 }
    @{ Sheet = "CRP"; Cell = "H89"; Value = 3082420 }
    @{ Sheet = "CRP"; Cell = "J89"; Value = 7059.4 }
    @{ Sheet = "CRP"; Cell = "L89"; Value = 35297 }
    @{ Sheet = "CRP"; Cell = "N89"; Value = -46529 }
    @{ Sheet = "CRP"; Cell = "H132"; Value = 3899.6216 }
    @{ Sheet = "CRP"; Cell = "I132"; Value = 3551.0286 }
    @{ Sheet = "CRP"; Cell = "K132"; Value = 10653.0858 }
    @{ Sheet = "CRP"; Cell = "M132"; Value = -8123.085800000001 }
    @{ Sheet = "CRP"; Cell = "H134"; Value = 2995.1516 }
    @{ Sheet = "CRP"; Cell = "I134"; Value = 2258.3333 }
    @{ Sheet = "CRP"; Cell = "K134"; Value = 6774.999899999999 }
    @{ Sheet = "CRP"; Cell = "M134"; Value = -4239.999899999999 }
    @{ Sheet = "CRP"; Cell = "H141"; Value = 300067.6 }
    @{ Sheet = "CRP"; Cell = "J141"; Value = 333747.62 }
    @{ Sheet = "CRP"; Cell = "L141"; Value = 333747.62 }
    @{ Sheet = "CRP"; Cell = "N141"; Value = -344107.62 }
    @{ Sheet = "CUL"; Cell = "H38"; Value = 13.25 }
    @{ Sheet = "CUL"; Cell = "J38"; Value = 21 }
    @{ Sheet = "CUL"; Cell = "L38"; Value = 63 }
    @{ Sheet = "CUL"; Cell = "N38"; Value = -757 }
    @{ Sheet = "CUL"; Cell = "H107"; Value = 444.24 }
    @{ Sheet = "CUL"; Cell = "I107"; Value = 281.25 }
    @{ Sheet = "CUL"; Cell = "J107"; Value = 475.2857 }
    @{ Sheet = "CUL"; Cell = "K107"; Value = 843.75 }
    @{ Sheet = "CUL"; Cell = "L107"; Value = 1425.8571 }
    @{ Sheet = "CUL"; Cell = "M107"; Value = 1076.25 }
    @{ Sheet = "CUL"; Cell = "N107"; Value = -5265.8571 }
    @{ Sheet = "CUL"; Cell = "H128"; Value = 236023.38 }
    @{ Sheet = "CUL"; Cell = "I128"; Value = 236023.38 }
    @{ Sheet = "CUL"; Cell = "K128"; Value = 708070.14 }
    @{ Sheet = "CUL"; Cell = "M128"; Value = -703090.14 }
    @{ Sheet = "CUL"; Cell = "H129"; Value = 12828130 }
    @{ Sheet = "CUL"; Cell = "J129"; Value = 15159562 }
    @{ Sheet = "CUL"; Cell = "L129"; Value = 45478686 }
    @{ Sheet = "CUL"; Cell = "N129"; Value = -45488686 }
    @{ Sheet = "CUL"; Cell = "H131"; Value = 4225.1377 }
    @{ Sheet = "CUL"; Cell = "J131"; Value = 4790.9546 }
    @{ Sheet = "CUL"; Cell = "L131"; Value = 14372.8638 }
    @{ Sheet = "CUL"; Cell = "N131"; Value = -24452.8638 }
    @{ Sheet = "CUL"; Cell = "H132"; Value = 2115.2144 }
    @{ Sheet = "CUL"; Cell = "J132"; Value = 2098.375 }
    @{ Sheet = "CUL"; Cell = "L132"; Value = 18885.375 }
    @{ Sheet = "CUL"; Cell = "N132"; Value = -23945.375 }
    @{ Sheet = "CUL"; Cell = "H136"; Value = 3152.2307 }
    @{ Sheet = "CUL"; Cell = "I136"; Value = 2397.2856 }
    @{ Sheet = "CUL"; Cell = "J136"; Value = 4033 }
    @{ Sheet = "CUL"; Cell = "K136"; Value = 7191.8568 }
    @{ Sheet = "CUL"; Cell = "L136"; Value = 12099 }
    @{ Sheet = "CUL"; Cell = "M136"; Value = -2091.8568 }
    @{ Sheet = "CUL"; Cell = "N136"; Value = -22299 }
    @{ Sheet = "CUL"; Cell = "H139"; Value = 2337.9583 }
    @{ Sheet = "CUL"; Cell = "I139"; Value = 1912.6818 }
    @{ Sheet = "CUL"; Cell = "J139"; Value = 7016 }
    @{ Sheet = "CUL"; Cell = "K139"; Value = 5738.0454 }
    @{ Sheet = "CUL"; Cell = "L139"; Value = 21048 }
    @{ Sheet = "CUL"; Cell = "M139"; Value = -598.0454 }
    @{ Sheet = "CUL"; Cell = "N139"; Value = -31328 }
    @{ Sheet = "CUL"; Cell = "H140"; Value = 1542.4722 }
    @{ Sheet = "CUL"; Cell = "I140"; Value = 645.9231 }
    @{ Sheet = "CUL"; Cell = "J140"; Value = 2049.2173 }
    @{ Sheet = "CUL"; Cell = "K140"; Value = 1937.7693 }
    @{ Sheet = "CUL"; Cell = "L140"; Value = 6147.651899999999 }
    @{ Sheet = "CUL"; Cell = "M140"; Value = 3242.2307 }
    @{ Sheet = "CUL"; Cell = "N140"; Value = -16507.6519 }
    @{ Sheet = "GSM"; Cell = "H80"; Value = 9000 }
    @{ Sheet = "GSM"; Cell = "I80"; Value = 8000 }
    @{ Sheet = "GSM"; Cell = "J80"; Value = 10000 }
    @{ Sheet = "GSM"; Cell = "K80"; Value = 8000 }
    @{ Sheet = "GSM"; Cell = "L80"; Value = 10000 }
    @{ Sheet = "GSM"; Cell = "M80"; Value = -7002 }
    @{ Sheet = "GSM"; Cell = "N80"; Value = -11996 }
    @{ Sheet = "GSM"; Cell = "H83"; Value = 9000 }
    @{ Sheet = "GSM"; Cell = "I83"; Value = 8000 }
    @{ Sheet = "GSM"; Cell = "J83"; Value = 10000 }
    @{ Sheet = "GSM"; Cell = "K83"; Value = 40000 }
    @{ Sheet = "GSM"; Cell = "L83"; Value = 50000 }
    @{ Sheet = "GSM"; Cell = "M83"; Value = -35008 }
    @{ Sheet = "GSM"; Cell = "N83"; Value = -59984 }
    @{ Sheet = "GSM"; Cell = "H132"; Value = 8688.454 }
    @{ Sheet = "GSM"; Cell = "J132"; Value = 11859.25 }
    @{ Sheet = "GSM"; Cell = "L132"; Value = 35577.75 }
    @{ Sheet = "GSM"; Cell = "N132"; Value = -40637.75 }
    @{ Sheet = "GSM"; Cell = "H136"; Value = 29609.826 }
    @{ Sheet = "GSM"; Cell = "J136"; Value = 29609.826 }
    @{ Sheet = "GSM"; Cell = "L136"; Value = 88829.478 }
    @{ Sheet = "GSM"; Cell = "N136"; Value = -93929.478 }
    @{ Sheet = "LTW"; Cell = "H16"; Value = 2413.625 }
    @{ Sheet = "LTW"; Cell = "J16"; Value = 2482.8333 }
    @{ Sheet = "LTW"; Cell = "L16"; Value = 2482.8333 }
    @{ Sheet = "LTW"; Cell = "N16"; Value = -2822.8333 }
    @{ Sheet = "LTW"; Cell = "H22"; Value = 4145.857 }
    @{ Sheet = "LTW"; Cell = "I22"; Value = 5156 }
    @{ Sheet = "LTW"; Cell = "J22"; Value = 2799 }
    @{ Sheet = "LTW"; Cell = "K22"; Value = 5156 }
    @{ Sheet = "LTW"; Cell = "L22"; Value = 2799 }
    @{ Sheet = "LTW"; Cell = "M22"; Value = -4861 }
    @{ Sheet = "LTW"; Cell = "N22"; Value = -3389 }
    @{ Sheet = "LTW"; Cell = "H27"; Value = 4145.857 }
    @{ Sheet = "LTW"; Cell = "I27"; Value = 5156 }
    @{ Sheet = "LTW"; Cell = "J27"; Value = 2799 }
    @{ Sheet = "LTW"; Cell = "K27"; Value = 5156 }
    @{ Sheet = "LTW"; Cell = "L27"; Value = 2799 }
    @{ Sheet = "LTW"; Cell = "M27"; Value = -5049 }
    @{ Sheet = "LTW"; Cell = "N27"; Value = -3013 }
    @{ Sheet = "LTW"; Cell = "H40"; Value = 3528.2424 }
    @{ Sheet = "LTW"; Cell = "I40"; Value = 1990.5 }
    @{ Sheet = "LTW"; Cell = "K40"; Value = 1990.5 }
    @{ Sheet = "LTW"; Cell = "M40"; Value = -1854.5 }
    @{ Sheet = "LTW"; Cell = "H93"; Value = 2076.55 }
    @{ Sheet = "LTW"; Cell = "J93"; Value = 3196.375 }
    @{ Sheet = "LTW"; Cell = "L93"; Value = 3196.375 }
    @{ Sheet = "LTW"; Cell = "N93"; Value = -5692.375 }
    @{ Sheet = "LTW"; Cell = "H132"; Value = 4609.933 }
    @{ Sheet = "LTW"; Cell = "I132"; Value = 2457 }
    @{ Sheet = "LTW"; Cell = "K132"; Value = 7371 }
    @{ Sheet = "LTW"; Cell = "M132"; Value = -4841 }
    @{ Sheet = "WVR"; Cell = "H2"; Value = 57000 }
    @{ Sheet = "WVR"; Cell = "J2"; Value = 57000 }
    @{ Sheet = "WVR"; Cell = "L2"; Value = 57000 }
    @{ Sheet = "WVR"; Cell = "N2"; Value = -57224 }
    @{ Sheet = "WVR"; Cell = "H41"; Value = 17599 }
    @{ Sheet = "WVR"; Cell = "I41"; Value = 0 }
    @{ Sheet = "WVR"; Cell = "J41"; Value = 17599 }
    @{ Sheet = "WVR"; Cell = "K41"; Value = 0 }
    @{ Sheet = "WVR"; Cell = "L41"; Value = 17599 }
    @{ Sheet = "WVR"; Cell = "M41"; Value = $null }
    @{ Sheet = "WVR"; Cell = "N41"; Value = -18379 }
    @{ Sheet = "WVR"; Cell = "H117"; Value = 24999.5 }
    @{ Sheet = "WVR"; Cell = "J117"; Value = 24999.5 }
    @{ Sheet = "WVR"; Cell = "L117"; Value = 24999.5 }
    @{ Sheet = "WVR"; Cell = "N117"; Value = -34177.5 }
    @{ Sheet = "WVR"; Cell = "H132"; Value = 3042.95 }
    @{ Sheet = "WVR"; Cell = "I132"; Value = 2843.8823 }
    @{ Sheet = "WVR"; Cell = "J132"; Value = 4171 }
    @{ Sheet = "WVR"; Cell = "K132"; Value = 8531.6469 }
    @{ Sheet = "WVR"; Cell = "L132"; Value = 12513 }
    @{ Sheet = "WVR"; Cell = "M132"; Value = -6001.6469 }
    @{ Sheet = "WVR"; Cell = "N132"; Value = -17573 }
)

foreach ($edit in $edits) {
    $ws = $wb.Worksheets.Item($edit.Sheet)
    $ws.Range($edit.Cell).Value = $edit.Value
}

Write-Host "Applied $($edits.Count) cell updates."
